$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 10858.1

$ws.Range("B3").Value = 10731.4
$ws.Range("C3").Value = 10591

$ws.Range("C4").Value = 10057.9

$ws.Range("C5").Value = 9755.1

$ws.Range("C6").Value = 9774.799999999999

$ws.Range("C7").Value = 9705.5

$ws.Range("C8").Value = 9509.9

$ws.Range("C9").Value = 9204.200000000001

$ws.Range("C10").Value = 9760.200000000001

$ws.Range("C11").Value = 9396.1

$ws.Range("C12").Value = 9113.200000000001

$ws.Range("C14").Value = 9194.1

$ws.Range("C15").Value = 9682.299999999999

$ws.Range("C16").Value = 9391.5

$ws.Range("C18").Value = 9511.799999999999

$ws.Range("C19").Value = 9799

$ws.Range("C20").Value = 10246

$ws.Range("C21").Value = 9597.5

$ws.Range("C22").Value = 8791.4

$ws.Range("C23").Value = 8406.5

$ws.Range("C24").Value = 8433.799999999999
